$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.967.17'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -1.32%  '
$c.ClearFormats()

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.831.26'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '
$c.ClearFormats()

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.65%  '
$c.ClearFormats()

# Row 5
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -0.94%  '
$c.ClearFormats()

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.55%  '
$c.ClearFormats()

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4579'
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.92%  '
$c.ClearFormats()

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3712'
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.32%  '
$c.ClearFormats()

# Row 9
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.88%  '
$c.ClearFormats()

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8776'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.59%  '
$c.ClearFormats()

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07805'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.ClearFormats()

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '19.68'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -0.77%  '
$c.ClearFormats()

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.823.44'
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.68%  '
$c.ClearFormats()

# Row 14
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.79%  '
$c.ClearFormats()

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.408'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -1.91%  '
$c.ClearFormats()

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '87.33'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -5.03%  '
$c.ClearFormats()

# Row 17
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +0.72%  '
$c.ClearFormats()

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008722'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -1.63%  '
$c.ClearFormats()

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +0.52%  '
$c.ClearFormats()

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '26.991.88'
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -1.32%  '
$c.ClearFormats()

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '14.51'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.ClearFormats()

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.011'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -1.98%  '
$c.ClearFormats()

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.057.15'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -1.28%  '
$c.ClearFormats()

# Row 24
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c.ClearFormats()

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.028'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +7.46%  '
$c.ClearFormats()

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '151.43'
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.ClearFormats()

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.23'
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.70%  '
$c.ClearFormats()

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.971'
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -4.78%  '
$c.ClearFormats()

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '114.00'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -1.84%  '
$c.ClearFormats()

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.938'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -3.53%  '
$c.ClearFormats()

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08814'
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.40%  '
$c.ClearFormats()

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.030'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '
$c.ClearFormats()

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7532'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -1.73%  '
$c.ClearFormats()

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.483'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -0.19%  '
$c.ClearFormats()

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.138'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -2.70%  '
$c.ClearFormats()

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.570'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -1.94%  '
$c.ClearFormats()

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.091'
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c.ClearFormats()

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01945'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -0.98%  '
$c.ClearFormats()

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.05155'
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -1.18%  '
$c.ClearFormats()

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.892'
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -3.32%  '
$c.ClearFormats()

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.963'
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -0.92%  '
$c.ClearFormats()

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.4992'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -3.03%  '
$c.ClearFormats()

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1602'
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -2.21%  '
$c.ClearFormats()

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '8.336'
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c.ClearFormats()

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.4693'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -2.86%  '
$c.ClearFormats()

# Row 46
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +0.59%  '
$c.ClearFormats()

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.13'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -1.88%  '
$c.ClearFormats()

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '102.48'
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.66%  '
$c.ClearFormats()

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.615'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -2.18%  '
$c.ClearFormats()

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06119'
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.ClearFormats()

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '64.57'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -1.56%  '
$c.ClearFormats()
